# tc_02 changes pushed work in progress
# Adds a "subject" column (H) with header + one data value, shrinks the
# header row, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H: header "subject" --------------------------------------
# Matches the other header cells in the table (D1/E1): same fill/border
# combo, default font, no wrap.
$ws.Range("H1").Value = "subject"
$ws.Range("D1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# --- New column H: data "Software Engineering" ----------------------------
# Matches the small monospace font used in F2 (Consolas, 7pt) but with an
# explicit theme text color and without the wrapped alignment.
$ws.Range("H2").Value = "Software Engineering"
$ws.Range("F2").Copy()
$ws.Range("H2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H2").Font.ThemeColor = 1
$ws.Range("H2").WrapText = $false

# --- Header row shrinks now that the row no longer needs as much height --
$ws.Rows.Item(1).RowHeight = 28.8

# --- Move the active selection --------------------------------------------
$ws.Range("D12").Select()
